# Actualizar 03-10-2021 17-13-57
# Appends the new daily UF / IVP rows (803-833, dates 44265-44295) to the
# UF_IVP_DIARIO sheet, extends the named range / dimension accordingly,
# and widens columns B and C to fit the (no-longer-wrapped) header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UF_IVP_DIARIO")

# ------------------------------------------------------------------
# 1. New data rows (FECHA, Unidad de fomento (UF), Indice de valor
#    promedio (IVP)) for 2021-03-10 .. 2021-04-09.
# ------------------------------------------------------------------
$newRows = @(
    @(44265, 29355.01, 30451.39),
    @(44266, 29356.9, 30455.44),
    @(44267, 29358.8, 30459.48),
    @(44268, 29360.69, 30463.53),
    @(44269, 29362.58, 30467.58),
    @(44270, 29364.47, 30471.63),
    @(44271, 29366.37, 30475.68),
    @(44272, 29368.26, 30479.73),
    @(44273, 29370.15, 30483.79),
    @(44274, 29372.04, 30487.84),
    @(44275, 29373.94, 30491.89),
    @(44276, 29375.83, 30495.94),
    @(44277, 29377.72, 30500),
    @(44278, 29379.62, 30504.05),
    @(44279, 29381.51, 30508.11),
    @(44280, 29383.41, 30512.16),
    @(44281, 29385.3, 30516.22),
    @(44282, 29387.19, 30520.28),
    @(44283, 29389.09, 30524.33),
    @(44284, 29390.98, 30528.39),
    @(44285, 29392.88, 30532.45),
    @(44286, 29394.77, 30536.51),
    @(44287, 29396.67, 30540.57),
    @(44288, 29398.56, 30544.63),
    @(44289, 29400.45, 30548.69),
    @(44290, 29402.35, 30552.75),
    @(44291, 29404.24, 30556.81),
    @(44292, 29406.14, 30560.87),
    @(44293, 29408.04, 30564.93),
    @(44294, 29409.93, 30569),
    @(44295, 29411.83, 30573.06)
)

$firstNewRow = 803
$lastNewRow = $firstNewRow + $newRows.Count - 1

# Carry the number formatting / style of the last existing data row (802)
# down across the whole new block in one shot before filling in values.
$ws.Range("A802:C802").Copy() | Out-Null
$ws.Range(("A{0}:C{1}" -f $firstNewRow, $lastNewRow)).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$r = $firstNewRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 2. Extend the sheet-scoped defined name to the new bottom row.
# ------------------------------------------------------------------
$definedName = $wb.Names.Item("UF_IVP_DIARIO!UF_IVP_DIARIO")
$definedName.RefersTo = ("=UF_IVP_DIARIO!`$A`$1:`$C`${0}" -f $lastNewRow)

# ------------------------------------------------------------------
# 3. Header row (row 2) no longer needs its explicit wrapped height -
#    let Excel recompute it automatically (removes the ht="51" override).
# ------------------------------------------------------------------
$ws.Rows.Item(2).AutoFit() | Out-Null

# ------------------------------------------------------------------
# 4. Columns B and C grow to fit the (single-line) header text now that
#    it is no longer wrapped.
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 19
$ws.Columns.Item(3).ColumnWidth = 24

# ------------------------------------------------------------------
# 5. Bring the view down to the newly-added last row, mirroring the
#    author's final on-screen selection.
# ------------------------------------------------------------------
$ws.Range(("A{0}" -f $lastNewRow)).Select() | Out-Null
